$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells per diff
$ws.Range("G5").Value = "Delivered"
$ws.Range("E13").Value = "Paid"

# Add new row 15 (new order record)
$ws.Range("A15").Value = "'20250308004312"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = "2025-03-08 00:43:12"
$ws.Range("C15").Value = "[{'Item Name': 'Roti Sabzi', 'Price': 100.0}]"
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "Paid"
$ws.Range("F15").Value = "Take Away"
$ws.Range("G15").Value = "Delivered"
